$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SValimaki                                       ")

# Row 4: Inside entry
$ws.Range("A4").Value = "Inside"
$ws.Range("B4").Value = "SValimaki                                       "
$ws.Range("C4").Value = 43423
$ws.Range("D4").Value = 0.4584691899189815
$ws.Range("C2").Copy()
$ws.Range("C4").PasteSpecial(-4122)
$ws.Range("D2").Copy()
$ws.Range("D4").PasteSpecial(-4122)

# Row 5: Inside entry
$ws.Range("A5").Value = "Inside"
$ws.Range("B5").Value = "SValimaki                                       "
$ws.Range("C5").Value = 43423
$ws.Range("D5").Value = 0.4818042578125
$ws.Range("C2").Copy()
$ws.Range("C5").PasteSpecial(-4122)
$ws.Range("D2").Copy()
$ws.Range("D5").PasteSpecial(-4122)

# Row 6: Outside entry
$ws.Range("E6").Value = "Outside"
$ws.Range("F6").Value = "SValimaki                                       "
$ws.Range("G6").Value = 43423
$ws.Range("H6").Value = 0.4890109440046296
$ws.Range("G3").Copy()
$ws.Range("G6").PasteSpecial(-4122)
$ws.Range("H3").Copy()
$ws.Range("H6").PasteSpecial(-4122)

# Row 7: Inside entry
$ws.Range("A7").Value = "Inside"
$ws.Range("B7").Value = "SValimaki                                       "
$ws.Range("C7").Value = 43423
$ws.Range("D7").Value = 0.4937088486111111
$ws.Range("C2").Copy()
$ws.Range("C7").PasteSpecial(-4122)
$ws.Range("D2").Copy()
$ws.Range("D7").PasteSpecial(-4122)

$excel.CutCopyMode = 0
